$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6945182
$ws.Range("I43").Value = 749.75
$ws.Range("K43").Value = 749.75
$ws.Range("M43").Value = -680.75

$ws.Range("H55").Value = 237.5
$ws.Range("I55").Value = 206.25
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 206.25
$ws.Range("L55").Value = 300
$ws.Range("M55").Value = 7.75
$ws.Range("N55").Value = -728

$ws.Range("H100").Value = 2816.3333
$ws.Range("I100").Value = 2966.6667
$ws.Range("J100").Value = 2666
$ws.Range("K100").Value = 2966.6667
$ws.Range("L100").Value = 2666
$ws.Range("M100").Value = -2425.6667
$ws.Range("N100").Value = -3748

$ws.Range("H111").Value = 3064.65
$ws.Range("J111").Value = 4749.9
$ws.Range("L111").Value = 14249.7
$ws.Range("N111").Value = -20383.7

$ws.Range("H137").Value = 1567
$ws.Range("I137").Value = 1104.625
$ws.Range("J137").Value = 2800
$ws.Range("K137").Value = 3313.875
$ws.Range("L137").Value = 8400
$ws.Range("M137").Value = -763.875
$ws.Range("N137").Value = -13500

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3184.25
$ws.Range("I32").Value = 3064.889
$ws.Range("J32").Value = 4974.6665
$ws.Range("K32").Value = 3064.889
$ws.Range("L32").Value = 4974.6665
$ws.Range("M32").Value = -2777.889
$ws.Range("N32").Value = -5548.6665

$ws.Range("H74").Value = 1056.6
$ws.Range("I74").Value = 1069.5834
$ws.Range("J74").Value = 1004.6667
$ws.Range("K74").Value = 1069.5834
$ws.Range("L74").Value = 1004.6667
$ws.Range("M74").Value = -195.5834
$ws.Range("N74").Value = -2752.6667

$ws.Range("H77").Value = 1056.6
$ws.Range("I77").Value = 1069.5834
$ws.Range("J77").Value = 1004.6667
$ws.Range("K77").Value = 5347.916999999999
$ws.Range("L77").Value = 5023.3335
$ws.Range("M77").Value = -979.9169999999995
$ws.Range("N77").Value = -13759.3335

$ws.Range("H132").Value = 2912.8823
$ws.Range("I132").Value = 2710
$ws.Range("J132").Value = 3399.8
$ws.Range("K132").Value = 8130
$ws.Range("L132").Value = 10199.4
$ws.Range("M132").Value = -5600
$ws.Range("N132").Value = -15259.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1229.711
$ws.Range("I31").Value = 915.0454999999999
$ws.Range("J31").Value = 1530.6957
$ws.Range("K31").Value = 915.0454999999999
$ws.Range("L31").Value = 1530.6957
$ws.Range("M31").Value = -620.0454999999999
$ws.Range("N31").Value = -2120.6957

$ws.Range("H34").Value = 1229.711
$ws.Range("I34").Value = 915.0454999999999
$ws.Range("J34").Value = 1530.6957
$ws.Range("K34").Value = 915.0454999999999
$ws.Range("L34").Value = 1530.6957
$ws.Range("M34").Value = -713.0454999999999
$ws.Range("N34").Value = -1934.6957

$ws.Range("H132").Value = 3389.1633
$ws.Range("I132").Value = 3396.325
$ws.Range("J132").Value = 3357.3333
$ws.Range("K132").Value = 10188.975
$ws.Range("L132").Value = 10071.9999
$ws.Range("M132").Value = -7658.974999999999
$ws.Range("N132").Value = -15131.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 677940
$ws.Range("I4").Value = 700089.5
$ws.Range("J4").Value = 660009.4399999999
$ws.Range("K4").Value = 2100268.5
$ws.Range("L4").Value = 1980028.32
$ws.Range("M4").Value = -2100156.5
$ws.Range("N4").Value = -1980252.32

$ws.Range("H88").Value = 3566.6667
$ws.Range("J88").Value = 8000
$ws.Range("L88").Value = 24000
$ws.Range("N88").Value = -24856

$ws.Range("H91").Value = 3566.6667
$ws.Range("J91").Value = 8000
$ws.Range("L91").Value = 24000
$ws.Range("N91").Value = -26964

$ws.Range("H98").Value = 872.44446
$ws.Range("I98").Value = 310.33334
$ws.Range("J98").Value = 1996.6666
$ws.Range("K98").Value = 931.0000200000001
$ws.Range("L98").Value = 5989.9998
$ws.Range("M98").Value = 566.9999799999999
$ws.Range("N98").Value = -8985.9998

$ws.Range("H107").Value = 4344.074
$ws.Range("I107").Value = 655.0714
$ws.Range("J107").Value = 8316.846
$ws.Range("K107").Value = 1965.2142
$ws.Range("L107").Value = 24950.538
$ws.Range("M107").Value = -45.21420000000012
$ws.Range("N107").Value = -28790.538

$ws.Range("H113").Value = 605.9048
$ws.Range("I113").Value = 431.14285
$ws.Range("J113").Value = 693.2857
$ws.Range("K113").Value = 1293.42855
$ws.Range("L113").Value = 2079.8571
$ws.Range("M113").Value = 876.5714499999999
$ws.Range("N113").Value = -6419.8571

$ws.Range("H131").Value = 16950440
$ws.Range("J131").Value = 1375.3
$ws.Range("L131").Value = 4125.9
$ws.Range("N131").Value = -14205.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1653.125
$ws.Range("I102").Value = 1631
$ws.Range("J102").Value = 1706.8572
$ws.Range("K102").Value = 1631
$ws.Range("L102").Value = 1706.8572
$ws.Range("M102").Value = -9
$ws.Range("N102").Value = -4950.8572

$ws.Range("H113").Value = 1878.9286
$ws.Range("I113").Value = 1056.1111
$ws.Range("K113").Value = 1056.1111
$ws.Range("M113").Value = 1113.8889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1134.7368
$ws.Range("I16").Value = 876.1667
$ws.Range("J16").Value = 1578
$ws.Range("K16").Value = 876.1667
$ws.Range("L16").Value = 1578
$ws.Range("M16").Value = -706.1667
$ws.Range("N16").Value = -1918

$ws.Range("H53").Value = 5250
$ws.Range("I53").Value = 500
$ws.Range("K53").Value = 500
$ws.Range("M53").Value = 18

$ws.Range("H82").Value = 2537.1428
$ws.Range("I82").Value = 2692
$ws.Range("J82").Value = 2150
$ws.Range("K82").Value = 2692
$ws.Range("L82").Value = 2150
$ws.Range("M82").Value = -2331
$ws.Range("N82").Value = -2872

$ws.Range("H85").Value = 2537.1428
$ws.Range("I85").Value = 2692
$ws.Range("J85").Value = 2150
$ws.Range("K85").Value = 2692
$ws.Range("L85").Value = 2150
$ws.Range("M85").Value = -1444
$ws.Range("N85").Value = -4646

$ws.Range("H106").Value = 32600
$ws.Range("J106").Value = 32600
$ws.Range("L106").Value = 32600
$ws.Range("N106").Value = -35124

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16262

$ws.Range("H113").Value = 522.86957
$ws.Range("I113").Value = 273
$ws.Range("J113").Value = 1230.8334
$ws.Range("K113").Value = 819
$ws.Range("L113").Value = 3692.5002
$ws.Range("M113").Value = 1351
$ws.Range("N113").Value = -8032.5002

$ws.Range("H132").Value = 4954.727
$ws.Range("I132").Value = 6651.8335
$ws.Range("J132").Value = 2918.2
$ws.Range("K132").Value = 19955.5005
$ws.Range("L132").Value = 8754.599999999999
$ws.Range("M132").Value = -17425.5005
$ws.Range("N132").Value = -13814.6

$ws.Range("H141").Value = 44371.668
$ws.Range("J141").Value = 44371.668
$ws.Range("L141").Value = 44371.668
$ws.Range("N141").Value = -54731.668
